# Updated cryptos list on Thu Oct 10 14:46:52 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) cells whose new value would otherwise be auto-coerced into a
# Number by Excel's type inference (losing the exact text formatting, e.g.
# trailing zeros or decimal-exponent form). Force these to remain Text by
# pre-setting the cell's number format to Text ("@") before assigning, which
# mirrors formatting a cell as Text in the Excel UI prior to typing into it.
$textPriceCells = @(
  "D4","D5","D6","D12","D13","D14","D16","D19","D20","D21","D22","D25","D26",
  "D27","D28","D31","D34","D36","D37","D38","D39","D40","D41","D43","D44",
  "D46","D47","D48","D49","D50","D51"
)
foreach ($ref in $textPriceCells) {
  $ws.Range($ref).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "60.582.78"
$ws.Range("E2").Value = "  -1.95%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.406.65"
$ws.Range("E3").Value = "  -1.64%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "0.996"
$ws.Range("E4").Value = "  -0.25%  "

# Row 5 - BNB
$ws.Range("D5").Value = "564.44"
$ws.Range("E5").Value = "  -2.55%  "

# Row 6 - Solana
$ws.Range("D6").Value = "137.45"
$ws.Range("E6").Value = "  -2.53%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.15%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +1.10%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "2.389.32"
$ws.Range("E9").Value = "  -2.15%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -3.21%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -0.78%  "

# Row 12 - Toncoin
$ws.Range("D12").Value = "5.03"
$ws.Range("E12").Value = "  -2.55%  "

# Row 13 - Cardano
$ws.Range("D13").Value = "0.335"
$ws.Range("E13").Value = "  -1.13%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "25.63"
$ws.Range("E14").Value = "  -0.50%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "2.843.27"
$ws.Range("E15").Value = "  -1.61%  "

# Row 16 - ShibaInu
$ws.Range("D16").Value = "0.0000167"
$ws.Range("E16").Value = "  -2.56%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "60.653.70"
$ws.Range("E17").Value = "  -1.80%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.387.64"
$ws.Range("E18").Value = "  -2.15%  "

# Row 19 - Uniswap
$ws.Range("D19").Value = "8.02"
$ws.Range("E19").Value = "  +11.79%  "

# Row 20 - Chainlink
$ws.Range("D20").Value = "10.51"
$ws.Range("E20").Value = "  -0.77%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "322.13"
$ws.Range("E21").Value = "  -0.68%  "

# Row 22 - Polkadot
$ws.Range("D22").Value = "4.03"
$ws.Range("E22").Value = "  -0.92%  "

# Row 23 - LEO
$ws.Range("E23").Value = "  -8.02%  "

# Row 24 - Dai
$ws.Range("E24").Value = "  +0.07%  "

# Row 25 - SuiNetwork
$ws.Range("D25").Value = "1.80"
$ws.Range("E25").Value = "  -6.25%  "

# Row 26 - Litecoin
$ws.Range("D26").Value = "64.05"
$ws.Range("E26").Value = "  -1.33%  "

# Row 27 - Aptos
$ws.Range("D27").Value = "8.23"
$ws.Range("E27").Value = "  -9.43%  "

# Row 28 - Bittensor
$ws.Range("D28").Value = "549.77"
$ws.Range("E28").Value = "  -5.36%  "

# Row 29 - WrappedeETH
$ws.Range("D29").Value = "2.527.60"
$ws.Range("E29").Value = "  -1.37%  "

# Row 30 - PEPE
$ws.Range("D30").Value = "0.0₃0914"
$ws.Range("E30").Value = "  -0.87%  "

# Row 31 - InternetComputer(DFINITY)
$ws.Range("D31").Value = "7.90"
$ws.Range("E31").Value = "  +0.51%  "

# Row 32 - Fetch.AI
$ws.Range("E32").Value = "  -5.20%  "

# Row 33 - PancakeSwap
$ws.Range("E33").Value = "  -3.85%  "

# Row 34 - Kaspa
$ws.Range("D34").Value = "0.132"
$ws.Range("E34").Value = "  -1.38%  "

# Row 35 - FirstDigitalUSD
$ws.Range("E35").Value = "  +0.01%  "

# Row 36 / 37 swap: Monero <-> ImmutableX, with updated values
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "1.41"
$ws.Range("E36").Value = "  +1.14%  "

$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").Value = "152.67"
$ws.Range("E37").Value = "  +0.55%  "

# Row 38 - PolygonEcosystemToken
$ws.Range("D38").Value = "0.366"
$ws.Range("E38").Value = "  -1.32%  "

# Row 39 - NEARProtocol
$ws.Range("D39").Value = "4.52"
$ws.Range("E39").Value = "  -4.75%  "

# Row 40 - EthereumClassic
$ws.Range("D40").Value = "18.10"
$ws.Range("E40").Value = "  -0.77%  "

# Row 41 - RenderToken
$ws.Range("D41").Value = "5.07"
$ws.Range("E41").Value = "  -1.04%  "

# Row 42 - USDe
$ws.Range("E42").Value = "  +0.02%  "

# Row 43 - Stacks
$ws.Range("D43").Value = "1.64"
$ws.Range("E43").Value = "  -1.39%  "

# Row 44 - dogwifhat
$ws.Range("D44").Value = "2.34"
$ws.Range("E44").Value = "  +0.04%  "

# Row 45 - BabyDogeCoin
$ws.Range("D45").Value = "0.0₆0292"
$ws.Range("E45").Value = "  +5.77%  "

# Row 46 - Aave
$ws.Range("D46").Value = "141.92"
$ws.Range("E46").Value = "  +0.88%  "

# Row 47 - Filecoin
$ws.Range("D47").Value = "3.49"
$ws.Range("E47").Value = "  -1.76%  "

# Row 48 - Mantle
$ws.Range("D48").Value = "0.581"
$ws.Range("E48").Value = "  -2.54%  "

# Row 49 - Hedera
$ws.Range("D49").Value = "0.0498"
$ws.Range("E49").Value = "  -1.97%  "

# Row 50 - InjectiveProtocol
$ws.Range("D50").Value = "19.09"
$ws.Range("E50").Value = "  -2.34%  "

# Row 51 - Stellar
$ws.Range("D51").Value = "0.0895"
$ws.Range("E51").Value = "  -0.04%  "
